$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")

# --- 1) Create the new "2022-Q1" sheet ----------------------------------
# Clone the "2021-Q4" sheet (same column layout/headers/styles we need)
# and drop it in right before "总计", then rename + re-purpose it.
$srcSheet = $wb.Worksheets.Item("2021-Q4")
$srcSheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Copying a sheet in shifts everything after the insertion point, so the
# "总计" reference we grabbed earlier now points at the wrong sheet object
# (it was bound by position). Re-resolve it by name before using it again.
$totalSheet = $wb.Worksheets.Item("总计")

# The template sheet has two data rows; we only need one.
$newSheet.Rows.Item(3).Delete()

# Force the B:G columns to be stored as text (matches source data which
# keeps numeric-looking values such as "3.02" as plain text).
$newSheet.Range("B2:G2").NumberFormat = "@"

$newSheet.Range("B2").Value = "001703"
$newSheet.Range("C2").Value = "银华沪港深增长股票"
$newSheet.Range("D2").Value = "3.02"
$newSheet.Range("E2").Value = "87.52"
$newSheet.Range("F2").Value = "3.46"
$newSheet.Range("G2").Value = "0.1045"
$newSheet.Range("H2").Value = 9

# --- 2) Update the "总计" (totals) sheet --------------------------------
# Insert a new row for the 2022-Q1 totals above the existing 2021-Q4 row.
$totalSheet.Rows.Item(2).Insert()

# Re-use the format already sitting on A3 (the row that got pushed down,
# which still carries the original index-column style) for the new A2.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# The row-insert leaves stray formatting behind on B2:D2 - clear it so
# these cells look like their unstyled siblings in rows below.
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.1

# Renumber the index column for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
